$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.025178707422301
$ws.Cells.Item(2, 4).Value = 1.028335994001733
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.023673710081207
$ws.Cells.Item(2, 9).Value = 1.029360113836184
$ws.Cells.Item(2, 10).Value = 1.030349178630393
$ws.Cells.Item(2, 11).Value = 1.031153107958656
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.026504455809386
$ws.Cells.Item(2, 14).Value = 1.031812392646755
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.027084911950532
$ws.Cells.Item(3, 4).Value = 1.029735023772635
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.026231184010214
$ws.Cells.Item(3, 9).Value = 1.029839124794931
$ws.Cells.Item(3, 10).Value = 1.031889672209725
$ws.Cells.Item(3, 11).Value = 1.032358457520224
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.028864096660109
$ws.Cells.Item(3, 14).Value = 1.033355073903663
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.028311923916783
$ws.Cells.Item(4, 4).Value = 1.030634761780702
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.027878643086032
$ws.Cells.Item(4, 9).Value = 1.03014485212729
$ws.Cells.Item(4, 10).Value = 1.032879920949726
$ws.Cells.Item(4, 11).Value = 1.033132370767463
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.030383320860615
$ws.Cells.Item(4, 14).Value = 1.034346728910458
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028826251928956
$ws.Cells.Item(5, 4).Value = 1.031011710469309
$ws.Cells.Item(5, 5).Value = 0.9945786998344013
$ws.Cells.Item(5, 6).Value = 1.028569510991789
$ws.Cells.Item(5, 9).Value = 1.030272378200705
$ws.Cells.Item(5, 10).Value = 1.033294679312714
$ws.Cells.Item(5, 11).Value = 1.033456301753261
$ws.Cells.Item(5, 12).Value = 0.9971179600051301
$ws.Cells.Item(5, 13).Value = 1.031020222155952
$ws.Cells.Item(5, 14).Value = 1.034762076277895
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028912522354763
$ws.Cells.Item(6, 4).Value = 1.031074926151731
$ws.Cells.Item(6, 5).Value = 0.9946255319796335
$ws.Cells.Item(6, 6).Value = 1.028685411022077
$ws.Cells.Item(6, 9).Value = 1.030293731929852
$ws.Cells.Item(6, 10).Value = 1.033364229402389
$ws.Cells.Item(6, 11).Value = 1.033510608404331
$ws.Cells.Item(6, 12).Value = 0.9971555583673452
$ws.Cells.Item(6, 13).Value = 1.031127057541534
$ws.Cells.Item(6, 14).Value = 1.034831725136674
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.028318802275327
$ws.Cells.Item(7, 4).Value = 1.030639803670416
$ws.Cells.Item(7, 5).Value = 0.9943035907978917
$ws.Cells.Item(7, 6).Value = 1.027887881207598
$ws.Cells.Item(7, 9).Value = 1.030146560059822
$ws.Cells.Item(7, 10).Value = 1.03288546899523
$ws.Cells.Item(7, 11).Value = 1.033136704708876
$ws.Cells.Item(7, 12).Value = 0.9968970624459043
$ws.Cells.Item(7, 13).Value = 1.030391838104518
$ws.Cells.Item(7, 14).Value = 1.034352284834823
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.025824268609381
$ws.Cells.Item(8, 4).Value = 1.028809960699544
$ws.Cells.Item(8, 5).Value = 0.9929600610674297
$ws.Cells.Item(8, 6).Value = 1.024539580656781
$ws.Cells.Item(8, 9).Value = 1.029522878507822
$ws.Cells.Item(8, 10).Value = 1.030871168841028
$ws.Cells.Item(8, 11).Value = 1.031561722531735
$ws.Cells.Item(8, 12).Value = 0.9958175282591057
$ws.Cells.Item(8, 13).Value = 1.027303512180841
$ws.Cells.Item(8, 14).Value = 1.032335124143358
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.02137790619225
$ws.Cells.Item(9, 4).Value = 1.025542198134939
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.018580651290297
$ws.Cells.Item(9, 9).Value = 1.028391072200908
$ws.Cells.Item(9, 10).Value = 1.027270306464795
$ws.Cells.Item(9, 11).Value = 1.028739272303612
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.02180112033373
$ws.Cells.Item(9, 14).Value = 1.028729148129526
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.018377570427305
$ws.Cells.Item(10, 4).Value = 1.02333311824805
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.01456547641604
$ws.Cells.Item(10, 9).Value = 1.027613853700674
$ws.Cells.Item(10, 10).Value = 1.024833448134702
$ws.Cells.Item(10, 11).Value = 1.026824593080095
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.01808944770566
$ws.Cells.Item(10, 14).Value = 1.026288829181093
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.017069374125172
$ws.Cells.Item(11, 4).Value = 1.022368988112203
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.012816085277426
$ws.Cells.Item(11, 9).Value = 1.02727178742885
$ws.Cells.Item(11, 10).Value = 1.023769267090739
$ws.Cells.Item(11, 11).Value = 1.025987370742261
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.01647131839756
$ws.Cells.Item(11, 14).Value = 1.025223136877984
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.016582055176483
$ws.Cells.Item(12, 4).Value = 1.022009699670706
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.012164602482137
$ws.Cells.Item(12, 9).Value = 1.027143886159686
$ws.Cells.Item(12, 10).Value = 1.023372596055166
$ws.Cells.Item(12, 11).Value = 1.025675137305264
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.015868572459172
$ws.Cells.Item(12, 14).Value = 1.024825902524041
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.01668665059078
$ws.Cells.Item(13, 4).Value = 1.022086821554819
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.012304424789296
$ws.Cells.Item(13, 9).Value = 1.027171359749699
$ws.Cells.Item(13, 10).Value = 1.023457746645851
$ws.Cells.Item(13, 11).Value = 1.025742169519102
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.015997941369953
$ws.Cells.Item(13, 14).Value = 1.024911174038332
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.017029120911879
$ws.Cells.Item(14, 4).Value = 1.0223393132009
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.012762268146751
$ws.Cells.Item(14, 9).Value = 1.027261232332471
$ws.Cells.Item(14, 10).Value = 1.023736506628147
$ws.Cells.Item(14, 11).Value = 1.025961587121943
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.01642153024829
$ws.Cells.Item(14, 14).Value = 1.025190329891778
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.017239941913978
$ws.Cells.Item(15, 4).Value = 1.022494726068806
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.013044135999735
$ws.Cells.Item(15, 9).Value = 1.02731649376308
$ws.Cells.Item(15, 10).Value = 1.023908075062386
$ws.Cells.Item(15, 11).Value = 1.026096610864291
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.016682290302922
$ws.Cells.Item(15, 14).Value = 1.025362141972873
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.018464197455222
$ws.Cells.Item(16, 4).Value = 1.023396942212273
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.014681345288508
$ws.Cells.Item(16, 9).Value = 1.027636438036942
$ws.Cells.Item(16, 10).Value = 1.024903881718881
$ws.Cells.Item(16, 11).Value = 1.026879982707317
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.018196602193735
$ws.Cells.Item(16, 14).Value = 1.026359362789041
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.019229694893027
$ws.Cells.Item(17, 4).Value = 1.023960828405204
$ws.Cells.Item(17, 5).Value = 0.9894763578477731
$ws.Cells.Item(17, 6).Value = 1.015705392043278
$ws.Cells.Item(17, 9).Value = 1.027835642559039
$ws.Cells.Item(17, 10).Value = 1.02552609025161
$ws.Cells.Item(17, 11).Value = 1.027369169468966
$ws.Cells.Item(17, 12).Value = 0.9930127773692701
$ws.Cells.Item(17, 13).Value = 1.019143519449802
$ws.Cells.Item(17, 14).Value = 1.026982454929254
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.019675327950894
$ws.Cells.Item(18, 4).Value = 1.024289003870177
$ws.Cells.Item(18, 5).Value = 0.9897087662937551
$ws.Cells.Item(18, 6).Value = 1.016301663802155
$ws.Cells.Item(18, 9).Value = 1.027951302908891
$ws.Cells.Item(18, 10).Value = 1.02588814740184
$ws.Cells.Item(18, 11).Value = 1.02765371896383
$ws.Cells.Item(18, 12).Value = 0.9932001317071766
$ws.Cells.Item(18, 13).Value = 1.019694787495635
$ws.Cells.Item(18, 14).Value = 1.027345026242166
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.019827131161933
$ws.Cells.Item(19, 4).Value = 1.024400780278473
$ws.Cells.Item(19, 5).Value = 0.9897880325774039
$ws.Cells.Item(19, 6).Value = 1.016504802883105
$ws.Cells.Item(19, 9).Value = 1.027990650189327
$ws.Cells.Item(19, 10).Value = 1.026011453564924
$ws.Cells.Item(19, 11).Value = 1.027750610663994
$ws.Cells.Item(19, 12).Value = 0.993264023964098
$ws.Cells.Item(19, 13).Value = 1.019882578747501
$ws.Cells.Item(19, 14).Value = 1.027468507514145
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.019147654396246
$ws.Cells.Item(20, 4).Value = 1.023900404443342
$ws.Cells.Item(20, 5).Value = 0.9894336180355766
$ws.Cells.Item(20, 6).Value = 1.015595629245854
$ws.Cells.Item(20, 9).Value = 1.027814324944587
$ws.Cells.Item(20, 10).Value = 1.025459423040053
$ws.Cells.Item(20, 11).Value = 1.027316765733882
$ws.Cells.Item(20, 12).Value = 0.9929783193490043
$ws.Cells.Item(20, 13).Value = 1.019042033480619
$ws.Cells.Item(20, 14).Value = 1.02691569304261
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.016928310814271
$ws.Cells.Item(21, 4).Value = 1.022264993164338
$ws.Cells.Item(21, 5).Value = 0.9882828385668255
$ws.Cells.Item(21, 6).Value = 1.012627491576915
$ws.Cells.Item(21, 9).Value = 1.027234790458929
$ws.Cells.Item(21, 10).Value = 1.023654457319424
$ws.Cells.Item(21, 11).Value = 1.025897008879962
$ws.Cells.Item(21, 12).Value = 0.9920501090198107
$ws.Cells.Item(21, 13).Value = 1.016296841278972
$ws.Cells.Item(21, 14).Value = 1.025108164063625
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.015524820880342
$ws.Cells.Item(22, 4).Value = 1.021229975549505
$ws.Cells.Item(22, 5).Value = 0.9875604150241496
$ws.Cells.Item(22, 6).Value = 1.010751549384831
$ws.Cells.Item(22, 9).Value = 1.026865532584655
$ws.Cells.Item(22, 10).Value = 1.022511562788016
$ws.Cells.Item(22, 11).Value = 1.02499709529474
$ws.Cells.Item(22, 12).Value = 0.991467000034148
$ws.Cells.Item(22, 13).Value = 1.014560962161011
$ws.Cells.Item(22, 14).Value = 1.023963646490891
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.016269618770326
$ws.Cells.Item(23, 4).Value = 1.021779309208606
$ws.Cells.Item(23, 5).Value = 0.9879432794636459
$ws.Cells.Item(23, 6).Value = 1.011746966749919
$ws.Cells.Item(23, 9).Value = 1.027061750118132
$ws.Cells.Item(23, 10).Value = 1.023118206476667
$ws.Cells.Item(23, 11).Value = 1.025474853396
$ws.Cells.Item(23, 12).Value = 0.9917760702887607
$ws.Cells.Item(23, 13).Value = 1.015482138693955
$ws.Cells.Item(23, 14).Value = 1.024571151683162
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.019184727647765
$ws.Cells.Item(24, 4).Value = 1.023927709684762
$ws.Cells.Item(24, 5).Value = 0.9894529299347241
$ws.Cells.Item(24, 6).Value = 1.015645229532068
$ws.Cells.Item(24, 9).Value = 1.027823959101168
$ws.Cells.Item(24, 10).Value = 1.025489549759512
$ws.Cells.Item(24, 11).Value = 1.027340447150182
$ws.Cells.Item(24, 12).Value = 0.9929938892766438
$ws.Cells.Item(24, 13).Value = 1.019087893872143
$ws.Cells.Item(24, 14).Value = 1.026945862545465
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.022533608685352
$ws.Cells.Item(25, 4).Value = 1.026392275118564
$ws.Cells.Item(25, 5).Value = 0.9912096547607046
$ws.Cells.Item(25, 6).Value = 1.020128452817484
$ws.Cells.Item(25, 9).Value = 1.028687622087092
$ws.Cells.Item(25, 10).Value = 1.028207482838032
$ws.Cells.Item(25, 11).Value = 1.029474666719832
$ws.Cells.Item(25, 12).Value = 0.9944092447426411
$ws.Cells.Item(25, 13).Value = 1.023231061759744
$ws.Cells.Item(25, 14).Value = 1.029667655400708
